$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.2015706806282722
$ws.Range("C2").Value = 0.5340314136125655
$ws.Range("J2").Value = 0.01047120418848168
$ws.Range("P2").Value = 0.1596858638743456
$ws.Range("S2").Value = 0.09424083769633508

# Row 3
$ws.Range("B3").Value = 0.009615384615384616
$ws.Range("C3").Value = 0.009615384615384616
$ws.Range("J3").Value = 0.03365384615384615
$ws.Range("P3").Value = 0.7307692307692307
$ws.Range("S3").Value = 0.2163461538461539

# Row 6
$ws.Range("B6").Value = 0.08737864077669903
$ws.Range("D6").Value = 0.02427184466019417
$ws.Range("F6").Value = 0.04368932038834952
$ws.Range("J6").Value = 0.2184466019417476
$ws.Range("O6").Value = 0.02912621359223301
$ws.Range("Q6").Value = 0.1262135922330097
$ws.Range("R6").Value = 0.07766990291262135
$ws.Range("S6").Value = 0.3932038834951456

# Row 7
$ws.Range("B7").Value = 0.1515151515151515
$ws.Range("D7").Value = 0.01731601731601732
$ws.Range("F7").Value = 0.0303030303030303
$ws.Range("J7").Value = 0.1125541125541126
$ws.Range("O7").Value = 0.01731601731601732
$ws.Range("Q7").Value = 0.1645021645021645
$ws.Range("R7").Value = 0.04761904761904762
$ws.Range("S7").Value = 0.4588744588744589

# Row 8
$ws.Range("B8").Value = 0.1310592459605027
$ws.Range("D8").Value = 0.0125673249551167
$ws.Range("E8").Value = 0.001795332136445242
$ws.Range("F8").Value = 0.06822262118491922
$ws.Range("J8").Value = 0.1005385996409336
$ws.Range("O8").Value = 0.02154398563734291
$ws.Range("Q8").Value = 0.1813285457809695
$ws.Range("R8").Value = 0.059245960502693
$ws.Range("S8").Value = 0.4236983842010772

# Row 9
$ws.Range("B9").Value = 0.08333333333333333
$ws.Range("D9").Value = 0.01111111111111111
$ws.Range("F9").Value = 0.03888888888888889
$ws.Range("J9").Value = 0.1055555555555556
$ws.Range("O9").Value = 0.02222222222222222
$ws.Range("Q9").Value = 0.1722222222222222
$ws.Range("R9").Value = 0.05
$ws.Range("S9").Value = 0.5166666666666667

# Row 10
$ws.Range("B10").Value = 0.1232339089481947
$ws.Range("D10").Value = 0.0282574568288854
$ws.Range("F10").Value = 0.06122448979591837
$ws.Range("J10").Value = 0.1067503924646782
$ws.Range("O10").Value = 0.02354788069073784
$ws.Range("Q10").Value = 0.2150706436420722
$ws.Range("R10").Value = 0.0565149136577708
$ws.Range("S10").Value = 0.3854003139717426

# Row 11
$ws.Range("G11").Value = 0.1201201201201201
$ws.Range("J11").Value = 0.06606606606606606
$ws.Range("K11").Value = 0.1921921921921922
$ws.Range("L11").Value = 0.5945945945945946
$ws.Range("S11").Value = 0.02702702702702703

# Row 12
$ws.Range("G12").Value = 0.7549019607843137
$ws.Range("J12").Value = 0.1715686274509804
$ws.Range("L12").Value = 0.009803921568627451
$ws.Range("S12").Value = 0.06372549019607843

# Row 13
$ws.Range("F13").Value = 0.01724137931034483
$ws.Range("G13").Value = 0.7586206896551724
$ws.Range("J13").Value = 0.1896551724137931
$ws.Range("S13").Value = 0.03448275862068965

# Row 15
$ws.Range("F15").Value = 0.01333333333333333
$ws.Range("H15").Value = 0.1422222222222222
$ws.Range("I15").Value = 0.06666666666666667
$ws.Range("J15").Value = 0.3111111111111111
$ws.Range("K15").Value = 0.07555555555555556
$ws.Range("M15").Value = 0.008888888888888889
$ws.Range("O15").Value = 0.06666666666666667
$ws.Range("S15").Value = 0.3155555555555555

# Row 16
$ws.Range("F16").Value = 0.01619433198380567
$ws.Range("H16").Value = 0.1659919028340081
$ws.Range("I16").Value = 0.07692307692307693
$ws.Range("J16").Value = 0.3967611336032389
$ws.Range("K16").Value = 0.1133603238866397
$ws.Range("M16").Value = 0.03238866396761134
$ws.Range("N16").Value = 0.004048582995951417
$ws.Range("O16").Value = 0.02834008097165992
$ws.Range("S16").Value = 0.1659919028340081

# Row 17
$ws.Range("F17").Value = 0.01492537313432836
$ws.Range("H17").Value = 0.2196162046908316
$ws.Range("I17").Value = 0.08528784648187633
$ws.Range("J17").Value = 0.3816631130063966
$ws.Range("K17").Value = 0.07889125799573561
$ws.Range("M17").Value = 0.01918976545842218
$ws.Range("O17").Value = 0.07249466950959488
$ws.Range("S17").Value = 0.1279317697228145

# Row 18
$ws.Range("F18").Value = 0.02857142857142857
$ws.Range("H18").Value = 0.1857142857142857
$ws.Range("I18").Value = 0.1071428571428571
$ws.Range("J18").Value = 0.4
$ws.Range("K18").Value = 0.07142857142857142
$ws.Range("M18").Value = 0.01428571428571429
$ws.Range("O18").Value = 0.06428571428571428
$ws.Range("S18").Value = 0.1285714285714286

# Row 19
$ws.Range("F19").Value = 0.01557210561949898
$ws.Range("H19").Value = 0.2376438727149628
$ws.Range("I19").Value = 0.06431956668923494
$ws.Range("J19").Value = 0.3595125253893026
$ws.Range("K19").Value = 0.1157752200406229
$ws.Range("M19").Value = 0.02640487474610697
$ws.Range("O19").Value = 0.05145565335138795
$ws.Range("S19").Value = 0.1293161814488829

